$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6 and 7 swap their species-observation data (columns A, E, F, G, H, I),
# while column B ("Taxonsorteringsordning") receives new, independent values.
# Column D ("NT") is unchanged in both rows.

# Row 6 -> former row 7 species data (Spillkråka / Dryocopus martius), new B value
$ws.Range("A6").Value = 130889857
$ws.Range("B6").Value = 57881
$ws.Range("E6").Value = 100049
$ws.Range("F6").Value = "Spillkråka"
$ws.Range("G6").Value = "Dryocopus martius"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("I6").Value = "'1"
$ws.Range("I6").Style = "Normal"

# Row 7 -> former row 6 species data (Talltita / Poecile montanus), new B value
$ws.Range("A7").Value = 130889854
$ws.Range("B7").Value = 58043
$ws.Range("E7").Value = 103021
$ws.Range("F7").Value = "Talltita"
$ws.Range("G7").Value = "Poecile montanus"
$ws.Range("H7").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I7").Value = "'2"
$ws.Range("I7").Style = "Normal"
